$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CBSLine_BurstIncrease")

# Row 35: fill in MFS=176Byte case (mirrors row 4's AC figures)
$ws.Range("S35").Value = 1328.70144
$ws.Range("T35").Value = "5,632MBit/s"

# S35 picks up the same (non-bold, centered) number format as the
# analogous S4 cell once it holds a real value
$ws.Range("S4").Copy()
$ws.Range("S35").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 41 (TFA unshaped)
$ws.Range("N41").Value = 0.0002051340288
$ws.Range("O41").Value = 0.00022524702719999999
$ws.Range("P41").Value = 0.00024536002560000001
$ws.Range("Q41").Value = 0.00026547302400000002
$ws.Range("R41").Value = 0.00028558602239999999
$ws.Range("S41").Value = 0.00025049902080000002
$ws.Range("T41").Value = 0.0014772991487999999

# Row 42 (TFA link shaped)
$ws.Range("N42").Value = 0.00019264
$ws.Range("O42").Value = 0.00020672
$ws.Range("P42").Value = 0.00020672
$ws.Range("Q42").Value = 0.00020672
$ws.Range("R42").Value = 0.00020672
$ws.Range("S42").Value = 0.00015152000000000001
$ws.Range("T42").Value = 0.0011710399999999999

# Row 43 (TFA CBS shaped)
$ws.Range("N43").Value = 0.0002051340288
$ws.Range("O43").Value = 0.00022524702719999999
$ws.Range("P43").Value = 0.00024536002560000001
$ws.Range("Q43").Value = 0.00026547302400000002
$ws.Range("R43").Value = 0.00028558602239999999
$ws.Range("S43").Value = 0.00025049902080000002
$ws.Range("T43").Value = 0.0014772991487999999

# Row 44 (TFA combined shaped)
$ws.Range("N44").Value = 0.00019264
$ws.Range("O44").Value = 0.00020672
$ws.Range("P44").Value = 0.00020672
$ws.Range("Q44").Value = 0.00020672
$ws.Range("R44").Value = 0.00020672
$ws.Range("S44").Value = 0.00015152000000000001
$ws.Range("T44").Value = 0.0011710399999999999

# Row 45 (SFA)
$ws.Range("T45").Value = 0.0010427340288

# Row 46 (PMOO)
$ws.Range("T46").Value = 0.0010427340288

# Row 47 (TMA)
$ws.Range("T47").Value = 0.0010427340288

# Row 55 (TFA unshaped)
$ws.Range("N55").Value = 0.00021062446080000001
$ws.Range("O55").Value = 0.00032178662400000001
$ws.Range("T55").Value = 0.00053241108480000003

# Row 56 (TFA link shaped)
$ws.Range("N56").Value = 0.00017967999999999999
$ws.Range("O56").Value = 0.00017967999999999999
$ws.Range("T56").Value = 0.00035935999999999997

# Row 57 (TFA CBS shaped)
$ws.Range("N57").Value = 0.00021062446080000001
$ws.Range("O57").Value = 0.00032178662400000001
$ws.Range("T57").Value = 0.00053241108480000003

# Row 58 (TFA combined shaped)
$ws.Range("N58").Value = 0.00017967999999999999
$ws.Range("O58").Value = 0.00017967999999999999
$ws.Range("T58").Value = 0.00035935999999999997

# Row 59 (SFA)
$ws.Range("T59").Value = 0.0003339844608

# Row 60 (PMOO)
$ws.Range("T60").Value = 0.0003339844608

# Row 61 (TMA)
$ws.Range("T61").Value = 0.0003339844608

# Update the selection to match the saved workbook state
$ws.Range("X34").Select()
